# Settings.xlsx bugfix: "Fixed bugs with minor arrays"
#   - mega_arrays (G4) default value: True -> False
#   - Repeats (I4) default value: 200 -> 1
#   - selection moves from J13 to G4 (and the view scrolls so column D is
#     the left-most visible column)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- I4: Repeats 200 -> 1 ------------------------------------------------
# A plain numeric-looking string is kept as a shared string (text) cell by
# the engine, exactly like the original "200" value was, so a direct
# assignment is enough here.
$ws.Range("I4").Value2 = "1"

# --- G4: mega_arrays True -> False --------------------------------------
# G4 is a *text* cell (t="s") holding the literal word "True"/"False", not
# a real boolean. Assigning the literal string "False" straight to
# .Value2/.Value gets auto-typed into an actual Boolean by the engine
# (same as typing False into a cell in real Excel), which would change the
# cell's stored type. To keep it a text cell we compute the replacement
# text with a formula (so it is produced as a string-formula result, not
# "typed" input) in a scratch cell, then copy/paste-special just the
# resulting value on top of G4 - this preserves G4's existing style and
# writes a plain shared-string text cell.
$scratch = $ws.Range("ZZ1")
$scratch.Formula = '=SUBSTITUTE(G3,"True","False")'
$scratch.Copy()
$ws.Range("G4").PasteSpecial(-4163) # xlPasteValues
$scratch.Clear()
$excel.CutCopyMode = $false

# --- Selection / view ----------------------------------------------------
$ws.Range("G4").Select()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
